$wb = $excel.ActiveWorkbook

$sheet1 = $wb.Worksheets.Item("展览")
$sheet1.Range("F4").Value = 510
$sheet1.Range("F5").Value = 2372
$sheet1.Range("F6").Value = 12
$sheet1.Range("F7").Value = 73
$sheet1.Range("F8").Value = 81
$sheet1.Range("F9").Value = 1686
$sheet1.Range("F10").Value = 1686
$sheet1.Range("F12").Value = 77
$sheet1.Range("F13").Value = 1437
$sheet1.Range("F16").Value = 826
$sheet1.Range("F17").Value = 61
$sheet1.Range("F18").Value = 187
$sheet1.Range("F19").Value = 155
$sheet1.Range("F20").Value = 7475
$sheet1.Range("F21").Value = 8429
$sheet1.Range("F24").Value = 416
$sheet1.Range("F26").Value = 98
$sheet1.Range("F34").Value = 1497
$sheet1.Range("F38").Value = 302
$sheet1.Range("F39").Value = 32
$sheet1.Range("F40").Value = 784
$sheet1.Range("F44").Value = 268
$sheet1.Range("F45").Value = 218
$sheet1.Range("F46").Value = 97
$sheet1.Range("F49").Value = 29

$sheet3 = $wb.Worksheets.Item("本地生活")
$sheet3.Range("F3").Value = 2651
$sheet3.Range("F4").Value = 300
$sheet3.Range("F5").Value = 154

$sheet4 = $wb.Worksheets.Item("全部类型")
$sheet4.Range("F6").Value = 300
$sheet4.Range("F7").Value = 154
$sheet4.Range("F9").Value = 510
$sheet4.Range("F10").Value = 2372
$sheet4.Range("F11").Value = 12
$sheet4.Range("F12").Value = 73
$sheet4.Range("F13").Value = 81
$sheet4.Range("F14").Value = 1686
$sheet4.Range("F15").Value = 1686
$sheet4.Range("F17").Value = 1437
$sheet4.Range("F19").Value = 826
$sheet4.Range("F20").Value = 61
$sheet4.Range("F22").Value = 187
$sheet4.Range("F24").Value = 156
$sheet4.Range("F25").Value = 7475
$sheet4.Range("F26").Value = 8429
$sheet4.Range("F28").Value = 98
$sheet4.Range("F37").Value = 302
$sheet4.Range("F39").Value = 784
$sheet4.Range("F44").Value = 268
$sheet4.Range("F45").Value = 218
$sheet4.Range("F46").Value = 97
$sheet4.Range("F50").Value = 29
